$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("preguntas")

# Update header row: keep "item" (A1) and "pregunta" (B1) as-is,
# rename "escalas" -> "escala" (C1) and "posibles respuestas" -> "posibles_respuestas" (D1)
$ws.Range("C1").Value = "escala"
$ws.Range("D1").Value = "posibles_respuestas"
